$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text representation (e.g. "67.762.04",
# "0.998", "1.00") instead of being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.762.04'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.322.31'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.83'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.67'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.51%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.319.58'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.178'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.577'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.28'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000273'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '636.21'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +8.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.853.25'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.45'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.827.57'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.78%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.314.17'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.65'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.89'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.901'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.57'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.02'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.77'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.77'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.56'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.56'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.58'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.66'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '591.68'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.938.84'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.93'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.06%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.51'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.67%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.59'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.129'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.25'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.67'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '32.50'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.40'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0683'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.337'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0413'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.03%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +12.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.54'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.94'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.56%  '
